# Edit: Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# - Reorders/refreshes the "Periodo Mora" rows for JOHN JAIRO TEJERA MENDOZA
#   (now 13 periods, newest-first: 2507..2209) with updated Valor Mora /
#   Salario Basico figures, and appends a new worker (WILLIAM JOSE PACHECO
#   CASTANO) as an additional account-statement row.
# - Updates the summary block (Valor Mora total, worker count, period count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the two extra detail rows needed (12 -> 13 periods for
#    the existing worker, plus 1 new row for the new worker) by inserting
#    two rows right before the current last detail row (27), pushing it
#    down to row 29. This also carries row 27's "bottom of table" border
#    style down to the new last row, and shifts the two signature rows
#    (32/33 -> 34/35) plus their merged cells automatically.
# ---------------------------------------------------------------------
$ws.Rows("27:28").Insert()

# The freshly inserted rows (27:28) come back with blank default
# formatting, so restore the normal detail-row look (borders/number
# formats) by copying the format from row 26, a normal interior row.
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J28").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Rewrite the 13 detail rows (16-28) for JOHN JAIRO TEJERA MENDOZA,
#    newest period first, then the new row (29) for WILLIAM JOSE PACHECO
#    CASTANO.
# ---------------------------------------------------------------------
$periods = @("2507","2506","2505","2504","2503","2502","2501","2412","2411","2410","2211","2210","2209")
$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "9296780"
    $ws.Cells.Item($row, 4).Value = "JOHN JAIRO TEJERA MENDOZA"
    $ws.Cells.Item($row, 5).Value = $p
    if ($p -eq "2209") {
        $ws.Cells.Item($row, 6).Value = 25600
    } else {
        $ws.Cells.Item($row, 6).Value = 48000
    }
    $ws.Cells.Item($row, 7).Value = 828116
    $row = $row + 1
}

# New worker row (29)
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1143384728"
$ws.Range("D29").Value = "WILLIAM JOSE PACHECO CASTANO"
$ws.Range("E29").Value = "2505"
$ws.Range("F29").Value = 80000
$ws.Range("G29").Value = 2000000

# ---------------------------------------------------------------------
# 3) Update the summary block above the table.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 681600      # VALOR MORA (total)
$ws.Range("C13").Value = 2           # Cant. Trabajadores
$ws.Range("F13").Value = 13          # Cant. Periodos

# ---------------------------------------------------------------------
# 4) Column D needs to be a bit wider to fit the new, longer worker name.
# ---------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 31.5

$wb.Save()
